$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (Changed) date column C for rows 2-10 from 45184 to 45186
$ws.Range("C2:C10").Value = 45186
